$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (A2:BA2) gets re-spelled: the original 5 "word" cells (name, type,
# sources, transformations, RQ) are replaced by one character per cell of
# the literal text  ['name', 'type', 'sources', 'transformations', 'RQ' ]
# laid out across columns A..BA, all keeping the existing bold/centered
# header style (style index 1, same as A1 / old A2:E2).
# ---------------------------------------------------------------------------

# Make sure every cell in the run carries the header format (copy it from
# A2, which already uses that format) before filling in values.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A2:BA2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$headerStr = "['name', 'type', 'sources', 'transformations', 'RQ' ]"

# Scratch cell, far away from any used range, for round-tripping a literal
# apostrophe through a formula so it lands as plain text instead of
# triggering Excel's "quote prefix" (leading-apostrophe) interpretation.
$scratch = $ws.Cells.Item(1000, 1000)

for ($i = 0; $i -lt $headerStr.Length; $i++) {
    $ch = $headerStr.Substring($i, 1)
    $target = $ws.Cells.Item(2, $i + 1)
    if ($ch -eq "'") {
        $scratch.Formula = "=""'"""
        $scratch.Copy()
        $target.PasteSpecial(-4163) | Out-Null  # xlPasteValues
        $scratch.Clear()
    } else {
        $target.Value = $ch
    }
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# New data rows: an extra "env" column (CodeSpeedy) for both table rows.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "employee_id"
$ws.Range("B3").Value = "int64"
$ws.Range("C3").Value = "CodeSpeedy"

$ws.Range("A4").Value = "employee_name"
$ws.Range("B4").Value = "object"
$ws.Range("C4").Value = "CodeSpeedy"
